$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header-style formatting from H1 onto the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header text for new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for new columns I and J (rows 2-19)
$values = @{
    2  = 4
    3  = 6
    4  = 8
    5  = 6
    6  = 7
    7  = 9
    8  = 7
    9  = 8
    10 = 8
    11 = 7
    12 = 8
    13 = 9
    14 = 5
    15 = 6
    16 = 6
    17 = 5
    18 = 6
    19 = 5
}

foreach ($row in $values.Keys) {
    $val = $values[$row]
    $ws.Cells.Item($row, 9).Value = $val
    $ws.Cells.Item($row, 10).Value = $val
}
